# Generated the 1.6.2 bundle and bumped version to 1.6.3
#
# - "Roll Variants" (row 9) loses its second requestor (Edward Robbins);
#   that request is attributed to the new feature below instead.
# - "Edit Saved Rolls" (row 19) is marked completed in version 1.6.2.
# - New feature request row 24: "Add multiples of dice to roll",
#   requested by Edward Robbins - Store Review.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll Variants (row 9): drop the second requestor.
$ws.Range("E9").ClearContents()

# Edit Saved Rolls (row 19): completed in 1.6.2.
$ws.Range("C19").Value = "1.6.2"

# New feature request row.
$ws.Range("A24").Value = "Add multiples of dice to roll"
$ws.Range("B24").Value = "I want to roll 6 d20s each with their own modifier in one roll"
$ws.Range("D24").Value = "Edward Robbins - Store Review"

# Matches the cached selection left behind in the saved workbook.
$ws.Range("D28").Select()
